$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("F1").Value = "Job"
$ws.Range("G1").Value = "Location"
$ws.Range("H1").Value = "Phones"
$ws.Range("I1").Value = "Emails"

# --- Clear email values for rows without changes in the diff ---
$clearRows = @(2,3,4,5,7,9,11,12,13,14,15,19,20,21,22)
foreach ($r in $clearRows) {
    $ws.Cells.Item($r, 9).Value = ""
}

# --- Update email values that gained extra addresses ---
$ws.Range("I6").Value = "lorcan.mulvey@mcaleer-rushe.co.uk , lorcanmulvey@yahoo.ie , lorcan.mulvey@yahoo.ie , lorcan.mulvey@berkeleygroup.co.uk"
$ws.Range("I10").Value = "leergray3@hotmail.co.uk , lee.gray@mcaleer-rushe.co.uk"
$ws.Range("I16").Value = "connor.graham@patton.co.uk , connor.graham@mcaleer-rushe.co.uk"
$ws.Range("I17").Value = "cathal.magee@mcaleer-rushe.co.uk , cathal.magee1@hotmail.co.uk"
